# Update CaN_template_mini.xlsx: add "Active" column to the Constraints sheet
# and make the Constraints sheet the active one.

$wb = $excel.ActiveWorkbook

$wsConstraints = $wb.Worksheets.Item("Constraints")

# Add the "Active" header in column D and set a value of 1 (active) for
# every existing constraint row.
$wsConstraints.Range("D1").Value = "Active"
$wsConstraints.Range("D2").Value = 1
$wsConstraints.Range("D3").Value = 1
$wsConstraints.Range("D4").Value = 1

# Select cell D4 on the Constraints sheet and make it the active sheet/tab.
$wsConstraints.Range("D4").Select()
$wsConstraints.Activate()

# Restore simple single-cell selections on the other sheets (the original
# file carried leftover multi-range selections that included "H3:H4").
$wsInfo = $wb.Worksheets.Item("INFO")
$wsInfo.Range("D13").Select()

$wsComponents = $wb.Worksheets.Item("Components & input parameter")
$wsComponents.Range("H3").Select()

$wsFluxes = $wb.Worksheets.Item("Fluxes")
$wsFluxes.Range("B18").Select()

$wsInputTS = $wb.Worksheets.Item("Input time-series")
$wsInputTS.Range("A6").Select()

# Re-activate the Constraints sheet so it ends up as the active tab.
$wsConstraints.Activate()
$wsConstraints.Range("D4").Select()

$wb.Save()
